$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that should be removed entirely (E2, C3, E3, C4, E4, C5, E5, C6, E6)
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Update values for rows 7-19 in columns C and E
$ws.Range("C7").Value = 4.880442637054072
$ws.Range("E7").Value = 1.641301872652501

$ws.Range("C8").Value = 5.941867202078877
$ws.Range("E8").Value = 2.672847571394987

$ws.Range("C9").Value = 0.292749233164491
$ws.Range("E9").Value = 2.181874035977249

$ws.Range("C10").Value = 2.032207428223742
$ws.Range("E10").Value = 2.13692496326825

$ws.Range("C11").Value = 2.775332754349846
$ws.Range("E11").Value = 2.200426660963761

$ws.Range("C12").Value = 3.565025829754953
$ws.Range("E12").Value = 2.446228176258058

$ws.Range("C13").Value = 3.444206290325491
$ws.Range("E13").Value = 2.496958452261078

$ws.Range("C14").Value = 2.667234932970275
$ws.Range("E14").Value = 2.567662999186382

$ws.Range("C15").Value = -4.511102905979703
$ws.Range("E15").Value = 0.9582724917052587

$ws.Range("C16").Value = 1.386772772629241
$ws.Range("E16").Value = 0.8813242377093244

$ws.Range("C17").Value = -0.9537175292835154
$ws.Range("E17").Value = 0.9049225073274991

$ws.Range("C18").Value = -3.303819519576723
$ws.Range("E18").Value = -0.3655818470008065

$ws.Range("C19").Value = -1.704805397136089
$ws.Range("E19").Value = 0.1460701281005727
